{"js": "// Insert a new paragraph \"Change 02\" right after the first paragraph (the\n// Nigeria description) and before the image paragraph. The new paragraph\n// inherits the \"Normal (Web)\" / NormalWeb style automatically from the\n// paragraph it is split off from; we only need to (re)apply the direct\n// character formatting (Segoe UI, 10.5pt / 21 half-points) that the rest of\n// the document uses.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Insert the new paragraph immediately after the first (\"Nigeria is a West\n// African country...\") paragraph, so it lands between it and the image\n// paragraph that currently follows.\nconst newParagraph = firstParagraph.insertParagraph(\"Change 02\", Word.InsertLocation.after);\n\n// Match the direct run formatting used throughout the document.\nnewParagraph.font.nameAscii = \"Segoe UI\";\nnewParagraph.font.nameBidirectional = \"Segoe UI\";\nnewParagraph.font.size = 10.5; // 21 half-points\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph \"Change 02\" right after the first paragraph (the\n# Nigeria description) and before the image paragraph, matching the\n# NormalWeb / Segoe UI / 10.5pt (21 half-points) formatting used throughout\n# the document.\n\n$d = $word.ActiveDocument\n\n$firstParagraph = $d.Paragraphs.Item(1)\n$firstRange = $firstParagraph.Range\n\n# InsertParagraphAfter() splits in a new paragraph right after the first\n# one, inheriting its paragraph style / direct formatting.\n$firstRange.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Item(2)\n$newRange = $newParagraph.Range\n$newRange.Text = \"Change 02\"\n\n# Re-assert the direct run formatting used elsewhere in the document.\n$newRange.Font.NameAscii = \"Segoe UI\"\n$newRange.Font.NameBi = \"Segoe UI\"\n$newRange.Font.Size = 10.5\n"}
